$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Find the last used row in column C (the "Förändrad" date column)
$lastRow = $ws.Cells.Item($ws.Rows.Count, 3).End(-4162).Row  # xlUp = -4162

# Column C holds the "Förändrad" (last changed) date, stored as a date serial
# number. Every data row (2..lastRow) had this value bumped by exactly one
# day (45178 -> 45179) when the report was refreshed.
for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 3)
    $v = $cell.Value2
    if ($v -ne $null) {
        $cell.Value2 = $v + 1
    }
}
